$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. The previously-last sheet "总计" becomes the new quarter sheet "2022-Q1".
#    A brand-new "总计" sheet is appended after it, holding the refreshed
#    summary (with the new 2022-Q1 row on top).
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# ---------------------------------------------------------------------------
# 2. Populate the "2022-Q1" sheet (fund holdings detail, same layout as the
#    2021-Q3 / 2021-Q4 sheets).
# ---------------------------------------------------------------------------

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"
# B1:D1 already carry style 2 from the old sheet; extend it across E1:H1.
$q1.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# Row 2
$q1.Range("B2").Value = "'009394"
$q1.Range("C2").Value = "银华同力精选混合"
$q1.Range("D2").Value = "'20.03"
$q1.Range("E2").Value = "'94.68"
$q1.Range("F2").Value = "'5.19"
$q1.Range("G2").Value = "'1.0396"
$q1.Range("H2").Value = 8

# Row 3
$q1.Range("B3").Value = "'003751"
$q1.Range("C3").Value = "万家瑞隆混合"
$q1.Range("D3").Value = "'27.84"
$q1.Range("E3").Value = "'86.40"
$q1.Range("F3").Value = "'3.26"
$q1.Range("G3").Value = "'0.9076"
$q1.Range("H3").Value = 7

# Row 4 (new)
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'005106"
$q1.Range("C4").Value = "银华农业产业股票"
$q1.Range("D4").Value = "'13.24"
$q1.Range("E4").Value = "'93.41"
$q1.Range("F4").Value = "'4.66"
$q1.Range("G4").Value = "'0.6170"
$q1.Range("H4").Value = 8

# Row 5 (new)
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'180020"
$q1.Range("C5").Value = "银华成长先锋混合"
$q1.Range("D5").Value = "'3.05"
$q1.Range("E5").Value = "'79.81"
$q1.Range("F5").Value = "'5.46"
$q1.Range("G5").Value = "'0.1665"
$q1.Range("H5").Value = 5

# Column A on rows 2-3 already has style 2 from the old sheet; copy it to the
# new A4:A5 cells, and strip the stray quote-prefix style picked up by the
# text-forced cells above (re-apply the plain, unstyled format of a blank
# cell to D2:G5).
$q1.Range("A2").Copy()
$q1.Range("A4:A5").PasteSpecial(-4122)

$q1.Range("Z1").Copy()
$q1.Range("B2:G5").PasteSpecial(-4122)

$q1.Application.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Populate the brand-new "总计" sheet (quarter-over-quarter summary).
# ---------------------------------------------------------------------------

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 2.73

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 3
$total.Range("D3").Value = 0.68

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.43

# Apply the bold/bordered header style (style index 2 in the workbook) to
# the header row and to the index column, copying it from the sibling sheet
# so no new style entries are introduced.
$q1.Range("B1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$total.Range("A2:A4").PasteSpecial(-4122)

$total.Application.CutCopyMode = $false
